$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

# Fill in review-count / rank data for the hotel in row 2.
# Source data stores these as text (shared strings), so force text
# formatting before assigning, then restore the default style so no
# extra formatting is left behind on the cells.
$ws.Range("G2:I2").NumberFormat = "@"
$ws.Range("G2").Value = "4"
$ws.Range("H2").Value = "9"
$ws.Range("I2").Value = "4"
$ws.Range("G2:I2").Style = "Normal"
